$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-30 Tuesday", "2024-07-31 Wednesday"),
    @("73×47=3431", "54×38=2052"),
    @("54×11=594", "72×44=3168"),
    @("72×67=4824", "64×99=6336"),
    @("69×56=3864", "33×55=1815"),
    @("40×31=1240", "60×67=4020"),
    @("51×88=4488", "55×20=1100"),
    @("32×41=1312", "85×57=4845"),
    @("45×41=1845", "63×75=4725"),
    @("24×21=504", "32×87=2784"),
    @("49×17=833", "13×90=1170"),
    @("97×95=9215", "18×64=1152"),
    @("56×39=2184", "57×75=4275"),
    @("34×50=1700", "12×20=240"),
    @("83×15=1245", "42×53=2226"),
    @("47×79=3713", "38×47=1786"),
    @("34×29=986", "65×92=5980"),
    @("76×78=5928", "93×61=5673"),
    @("23×76=1748", "12×68=816"),
    @("30×46=1380", "89×16=1424"),
    @("99×92=9108", "22×91=2002"),
    @("56×85=4760", "55×79=4345"),
    @("74×20=1480", "16×56=896"),
    @("80×32=2560", "61×83=5063"),
    @("38×74=2812", "95×92=8740"),
    @("13×27=351", "76×38=2888")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
